# Updated stock price / fundamentals data rows and normalized the fixed_ticker
# column (and related shared strings) to all point to the NTAP ticker, and set
# shares_outstanding to a single consistent value, resolving issues caused by
# leftover rows from other tickers being merged into this file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row=2; D=28.83040191256564; E=29.47870635986328; F=30.43209268573849; G=27.3660019125833; H=199618386 }
    @{ Row=3; D=27.92943773403355; E=25.5924072265625; F=28.33554384091013; G=23.63849510876349; H=199618386 }
    @{ Row=4; D=24.02051053238257; E=24.62930679321289; F=24.88361582238044; G=22.25576908304182; H=199618386 }
    @{ Row=5; D=26.29854859553848; E=23.76402854919434; F=26.91861398377902; G=23.3454847817201; H=199618386 }
    @{ Row=6; D=16.93152434774502; E=19.39045906066895; F=19.95250081282472; G=16.30703202311279; H=199618386 }
    @{ Row=7; D=18.61860349645597; E=20.0647087097168; F=20.67773056163994; G=17.68335032614316; H=199618386 }
    @{ Row=8; D=20.82498242605244; E=27.38920783996582; F=28.10185105302864; G=20.44490705260912; H=199618386 }
    @{ Row=9; D=27.02183398671808; E=29.1077880859375; F=31.05042923168181; G=24.17156539552661; H=199618386 }
    @{ Row=10; D=30.856978921634; E=33.48242568969727; F=33.76257862750008; G=30.08055375682135; H=199618386 }
    @{ Row=11; D=32.17403722646799; E=32.56003189086914; F=33.61346633188219; G=31.41009724996153; H=199618386 }
    @{ Row=12; D=35.336034753297; E=31.24636459350586; F=35.42494109551426; G=30.5836121482868; H=199618386 }
    @{ Row=13; D=36.32382353963828; E=45.87973022460938; F=46.88647140627459; G=35.64995588389973; H=199618386 }
    @{ Row=14; D=50.00074544961775; E=49.33265686035156; F=50.81548793486878; G=42.36660923305432; H=199618386 }
    @{ Row=15; D=54.22461726060622; E=55.84302139282227; F=59.54572660171732; G=52.15666404686269; H=199618386 }
    @{ Row=16; D=63.61856666804504; E=71.32542419433594; F=71.84305311698134; G=63.5199691902756; H=199618386 }
    @{ Row=17; D=64.78337122562856; E=55.19959259033203; F=69.29872049002491; G=53.47434415653837; H=199618386 }
    @{ Row=18; D=53.29106407476171; E=54.18021774291992; F=56.87261362212891; G=50.69008338975419; H=199618386 }
    @{ Row=19; D=61.02393413270417; E=49.46757507324219; F=61.57543379965095; G=48.55676994346013; H=199618386 }
    @{ Row=20; D=49.18869937409067; E=40.47267532348633; F=50.26662225233951; G=37.51680406287193; H=199618386 }
    @{ Row=21; D=47.98642698110644; E=51.50571441650391; F=54.34494685712181; G=47.32337052378312; H=199618386 }
    @{ Row=22; D=45.83212568610659; E=40.02386856079102; F=52.22292224045783; G=38.45615129113301; H=199618386 }
    @{ Row=23; D=36.9114889840923; E=38.61963272094727; F=40.27575061694861; G=34.51835653467458; H=199618386 }
    @{ Row=24; D=38.72356501919104; E=41.5465202331543; F=41.64295708992768; G=35.1378904508386; H=199618386 }
    @{ Row=25; D=39.40195370464413; E=47.24512481689453; F=49.38980897752485; G=38.92338664409048; H=199618386 }
    @{ Row=26; D=59.66606974150145; E=55.88950729370117; F=63.9961660945531; G=54.60387125996161; H=199618386 }
    @{ Row=27; D=67.3005737041848; E=69.52925109863281; F=72.48584010793769; G=65.8627206381978; H=199618386 }
    @{ Row=28; D=72.33097004838527; E=80.40491485595703; F=82.5386790242259; G=70.56790541266341; H=199618386 }
    @{ Row=29; D=81.20788714330122; E=80.80784606933594; F=84.36274012674011; G=76.48924480143835; H=199618386 }
    @{ Row=30; D=79.23526011787162; E=71.64812469482422; F=84.46398396421768; G=65.85265564532629; H=199618386 }
    @{ Row=31; D=67.58323473613144; E=66.18502044677734; F=70.58202997163593; G=59.40554444659634; H=199618386 }
    @{ Row=32; D=66.01632413880586; E=66.87861633300781; F=73.33189688441495; G=65.20966849109313; H=199618386 }
    @{ Row=33; D=65.33424872192228; E=63.15768432617188; F=70.23852086685753; G=60.25248635206163; H=199618386 }
    @{ Row=34; D=62.18377990250935; E=60.78999328613281; F=65.68709266067646; G=59.62222116924617; H=199618386 }
    @{ Row=35; D=59.66110051383101; E=62.98351669311523; F=67.03686848752579; G=58.41756991854681; H=199618386 }
    @{ Row=36; D=74.05310307481521; E=73.28868103027344; F=76.4610201220047; G=70.84254215184997; H=199618386 }
    @{ Row=37; D=70.00772077155352; E=87.90884399414062; F=88.11084395999148; G=69.62295578664545; H=199618386 }
    @{ Row=38; D=84.21239274455721; E=86.22482299804688; F=87.99536744668526; G=81.07765044532302; H=199618386 }
    @{ Row=39; D=98.5298510608514; E=117.0673828125; F=118.0880679316807; G=97.44112076148321; H=199618386 }
    @{ Row=40; D=123.3924385529623; E=117.8196258544922; F=131.1416699030146; G=110.460778291313; H=199618386 }
    @{ Row=41; D=113.263559397181; E=120.2028656005859; F=132.7582995743523; G=112.9597221902534; H=199618386 }
    @{ Row=42; D=117.4167557067094; E=98.26736450195312; F=125.8050698651996; G=96.00291711140476; H=199618386 }
    @{ Row=43; D=89.17720953135699; E=98.25347137451172; F=101.2359481614489; G=88.33498184215667; H=199618386 }
    @{ Row=44; D=102.0660508656801; E=112.3015289306641; F=117.7279229814402; G=100.1244913860024; H=199618386 }
)

foreach ($r in $rowData) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D   # open_price
    $ws.Cells.Item($r.Row, 5).Value = $r.E   # close_price
    $ws.Cells.Item($r.Row, 6).Value = $r.F   # high_price
    $ws.Cells.Item($r.Row, 7).Value = $r.G   # low_price
    $ws.Cells.Item($r.Row, 8).Value = $r.H   # shares_outstanding
    $ws.Cells.Item($r.Row, 9).Value = "NTAP" # fixed_ticker
}
